$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (SerialNo 21): Tree - Complete binary tree impl
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Tree"
$ws.Range("C22").Value = "Complete binary tree impl"
$ws.Range("H22").Value = "Tree"

# Row 23 (SerialNo 22): Breadth First search - BFS in tree
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Breadth First search"
$ws.Range("C23").Value = "BFS in tree"
$ws.Range("H23").Value = "BFS"

# Row 24 (SerialNo 23): Threaded BT
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Threaded BT : Inorder traversal and conversion to Threaded BT"
$ws.Range("H24").Value = "ThreadedBT"

$ws.Range("H24").Select() | Out-Null
